$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 121, shifting all existing
# rows 121..187 down to 123..189 (dates/values carried along unchanged).
$ws.Rows.Item(121).Insert()
$ws.Rows.Item(121).Insert()

# Populate the first new row (121) with its final values.
$ws.Range("A121").Value = 9
$ws.Range("B121").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C121").Value = "Metropolitana"
$ws.Range("D121").Value = 44518
$ws.Range("E121").Value = 13
$ws.Range("F121").Value = 100112021
$ws.Range("G121").Value = "Ají"
$ws.Range("H121").Value = "Inferno"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 52
$ws.Range("K121").Value = 18000
$ws.Range("L121").Value = 20000
$ws.Range("M121").Value = 19000
$ws.Range("N121").Value = "$/caja 12 kilos"
$ws.Range("O121").Value = "Región de Arica y Parinacota"
$ws.Range("P121").Value = 1583
$ws.Range("Q121").Value = 12
$ws.Range("R121").Value = "Hortaliza"

# Populate the second new row (122) with its final values.
$ws.Range("A122").Value = 9
$ws.Range("B122").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C122").Value = "Metropolitana"
$ws.Range("D122").Value = 44518
$ws.Range("E122").Value = 13
$ws.Range("F122").Value = 100112021
$ws.Range("G122").Value = "Ají"
$ws.Range("H122").Value = "Inferno"
$ws.Range("I122").Value = "Segunda"
$ws.Range("J122").Value = 25
$ws.Range("K122").Value = 16000
$ws.Range("L122").Value = 16000
$ws.Range("M122").Value = 16000
$ws.Range("N122").Value = "$/caja 12 kilos"
$ws.Range("O122").Value = "Región de Arica y Parinacota"
$ws.Range("P122").Value = 1333
$ws.Range("Q122").Value = 12
$ws.Range("R122").Value = "Hortaliza"
